# Rename the embedded logo pictures' display names.
#
#   - Pearson Edexcel logo (in both the first-page footer and the
#     default/primary footer): "image2.png" -> "image1.png"
#   - BTEC logo (in both the first-page header and the default/primary
#     header): "image1.jpg" -> "image2.jpg"
#
# The alt-text/description on every one of these pictures uniquely
# identifies which logo it is, so we drive the rename off that instead of
# hard-coding which physical header/footer part holds which picture.

$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument

function Rename-LogoShapes($shapes) {
    foreach ($shp in $shapes) {
        $alt = $shp.AlternativeText
        if ($alt -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        } elseif ($alt -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
    }
}

foreach ($sec in $d.Sections) {
    Rename-LogoShapes $sec.Headers($wdHeaderFooterPrimary).Range.InlineShapes
    Rename-LogoShapes $sec.Headers($wdHeaderFooterFirstPage).Range.InlineShapes
    Rename-LogoShapes $sec.Footers($wdHeaderFooterPrimary).Range.InlineShapes
    Rename-LogoShapes $sec.Footers($wdHeaderFooterFirstPage).Range.InlineShapes
}
